$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'56.636.61"
$ws.Cells.Item(2, 5).Value = "'  +0.35%  "
$ws.Cells.Item(3, 4).Value = "'2.385.61"
$ws.Cells.Item(3, 5).Value = "'  -0.57%  "
$ws.Cells.Item(4, 5).Value = "'  -0.07%  "
$ws.Cells.Item(5, 4).Value = "'504.62"
$ws.Cells.Item(5, 5).Value = "'  +0.62%  "
$ws.Cells.Item(6, 4).Value = "'132.78"
$ws.Cells.Item(6, 5).Value = "'  +3.48%  "
$ws.Cells.Item(7, 5).Value = "'  -0.17%  "
$ws.Cells.Item(8, 4).Value = "'0.550"
$ws.Cells.Item(8, 5).Value = "'  -0.08%  "
$ws.Cells.Item(9, 4).Value = "'2.390.21"
$ws.Cells.Item(9, 5).Value = "'  -0.32%  "
$ws.Cells.Item(10, 5).Value = "'  +1.96%  "
$ws.Cells.Item(11, 5).Value = "'  +0.53%  "
$ws.Cells.Item(12, 4).Value = "'0.327"
$ws.Cells.Item(12, 5).Value = "'  +2.87%  "
$ws.Cells.Item(13, 4).Value = "'4.71"
$ws.Cells.Item(13, 5).Value = "'  +2.06%  "
$ws.Cells.Item(14, 5).Value = "'  -0.44%  "
$ws.Cells.Item(15, 4).Value = "'56.587.32"
$ws.Cells.Item(15, 5).Value = "'  -0.82%  "
$ws.Cells.Item(16, 4).Value = "'21.69"
$ws.Cells.Item(16, 5).Value = "'  +0.86%  "
$ws.Cells.Item(18, 2).Value = "'Chainlink"
$ws.Cells.Item(18, 3).Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).Value = "'10.17"
$ws.Cells.Item(18, 5).Value = "'  +0.49%  "
$ws.Cells.Item(19, 2).Value = "'Polkadot"
$ws.Cells.Item(19, 3).Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(19, 4).Value = "'4.05"
$ws.Cells.Item(19, 5).Value = "'  +0.74%  "
$ws.Cells.Item(20, 2).Value = "'BitcoinCash"
$ws.Cells.Item(20, 3).Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).Value = "'309.76"
$ws.Cells.Item(20, 5).Value = "'  +0.01%  "
$ws.Cells.Item(21, 2).Value = "'WrappedEther"
$ws.Cells.Item(21, 3).Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(21, 4).Value = "'1.911.79"
$ws.Cells.Item(21, 5).Value = "'  -19.48%  "
$ws.Cells.Item(22, 5).Value = "'  +1.59%  "
$ws.Cells.Item(23, 5).Value = "'  +0.00%  "
$ws.Cells.Item(24, 5).Value = "'  -3.56%  "
$ws.Cells.Item(25, 4).Value = "'66.20"
$ws.Cells.Item(25, 5).Value = "'  +1.13%  "
$ws.Cells.Item(26, 4).Value = "'0.998"
$ws.Cells.Item(26, 5).Value = "'  -0.43%  "
$ws.Cells.Item(27, 5).Value = "'  +0.21%  "
$ws.Cells.Item(28, 4).Value = "'0.368"
$ws.Cells.Item(28, 5).Value = "'  -1.34%  "
$ws.Cells.Item(29, 4).Value = "'7.35"
$ws.Cells.Item(30, 4).Value = "'175.47"
$ws.Cells.Item(30, 5).Value = "'  +0.53%  "
$ws.Cells.Item(31, 4).Value = "'0.0₃0725"
$ws.Cells.Item(31, 5).Value = "'  +2.05%  "
$ws.Cells.Item(32, 5).Value = "'  -0.52%  "
$ws.Cells.Item(33, 5).Value = "'  +2.45%  "
$ws.Cells.Item(34, 4).Value = "'5.87"
$ws.Cells.Item(34, 5).Value = "'  -3.70%  "
$ws.Cells.Item(35, 5).Value = "'  +0.13%  "
$ws.Cells.Item(36, 5).Value = "'  +0.04%  "
$ws.Cells.Item(37, 5).Value = "'  +0.16%  "
$ws.Cells.Item(38, 4).Value = "'1.19"
$ws.Cells.Item(38, 5).Value = "'  -0.43%  "
$ws.Cells.Item(39, 5).Value = "'  +1.39%  "
$ws.Cells.Item(40, 4).Value = "'36.77"
$ws.Cells.Item(40, 5).Value = "'  +2.60%  "
$ws.Cells.Item(41, 4).Value = "'0.817"
$ws.Cells.Item(41, 5).Value = "'  +6.44%  "
$ws.Cells.Item(42, 5).Value = "'  +0.90%  "
$ws.Cells.Item(43, 4).Value = "'132.53"
$ws.Cells.Item(43, 5).Value = "'  +2.25%  "
$ws.Cells.Item(44, 4).Value = "'3.37"
$ws.Cells.Item(45, 5).Value = "'  +0.48%  "
$ws.Cells.Item(46, 4).Value = "'0.567"
$ws.Cells.Item(46, 5).Value = "'  -0.76%  "
$ws.Cells.Item(47, 4).Value = "'0.0908"
$ws.Cells.Item(47, 5).Value = "'  +1.14%  "
$ws.Cells.Item(48, 4).Value = "'246.57"
$ws.Cells.Item(48, 5).Value = "'  -2.64%  "
$ws.Cells.Item(49, 4).Value = "'0.0484"
$ws.Cells.Item(49, 5).Value = "'  +0.58%  "
$ws.Cells.Item(50, 5).Value = "'  +1.40%  "
$ws.Cells.Item(51, 4).Value = "'17.14"
$ws.Cells.Item(51, 5).Value = "'  +7.30%  "
